$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Datos actualizados" timestamp string (A1) ---
$ws.Range("A1").Value2 = "Datos actualizados a 27 de Marzo de 2020 a las 02:12"

# --- Reorder country rows: Camboya/Costa de Marfil swap (rows 101-102) ---
# --- and Honduras inserted ahead of Cuba, shifting rows 113-118 ---
$ws.Range("A101").Value2 = "Camboya"
$ws.Range("A102").Value2 = "Costa de Marfil"
$ws.Range("A113").Value2 = "Honduras"
$ws.Range("A114").Value2 = "Cuba"
$ws.Range("A115").Value2 = "Trinidad yTobago"
$ws.Range("A116").Value2 = "Nigeria"
$ws.Range("A117").Value2 = "Bolivia"
$ws.Range("A118").Value2 = "Liechtenstein"

# --- Update numeric statistics cells to the refreshed data snapshot ---
$ws.Range("B4").Value2 = 85344
$ws.Range("C4").Value2 = 17133
$ws.Range("D4").Value2 = 1868
$ws.Range("E4").Value2 = 82181
$ws.Range("G4").Value2 = 268
$ws.Range("H4").Value2 = 1295
$ws.Range("B5").Value2 = 81340
$ws.Range("C5").Value2 = 55
$ws.Range("D5").Value2 = 74588
$ws.Range("E5").Value2 = 3460
$ws.Range("F5").Value2 = 1034
$ws.Range("G5").Value2 = 5
$ws.Range("H5").Value2 = 3292
$ws.Range("B20").Value2 = 3370
$ws.Range("C20").Value2 = 286
$ws.Range("E20").Value2 = 3350
$ws.Range("E82").Value2 = 190
$ws.Range("G82").Value2 = 0
$ws.Range("H82").Value2 = 5
$ws.Range("B101").Value2 = 98
$ws.Range("C101").Value2 = 2
$ws.Range("D101").Value2 = 10
$ws.Range("E101").Value2 = 88
$ws.Range("F101").Value2 = 1
$ws.Range("C102").Value2 = 16
$ws.Range("D102").Value2 = 3
$ws.Range("E102").Value2 = 93
$ws.Range("F102").Value2 = 0
$ws.Range("C113").Value2 = 15
$ws.Range("D113").Value2 = 0
$ws.Range("E113").Value2 = 66
$ws.Range("F113").Value2 = 0
$ws.Range("H113").Value2 = 1
$ws.Range("B114").Value2 = 67
$ws.Range("C114").Value2 = 10
$ws.Range("D114").Value2 = 1
$ws.Range("F114").Value2 = 2
$ws.Range("G114").Value2 = 1
$ws.Range("H114").Value2 = 2
$ws.Range("C115").Value2 = 5
$ws.Range("D115").Value2 = 0
$ws.Range("E115").Value2 = 64
$ws.Range("B116").Value2 = 65
$ws.Range("C116").Value2 = 14
$ws.Range("D116").Value2 = 3
$ws.Range("H116").Value2 = 1
$ws.Range("B117").Value2 = 61
$ws.Range("C117").Value2 = 29
$ws.Range("E117").Value2 = 61
$ws.Range("B118").Value2 = 56
$ws.Range("C118").Value2 = 5
$ws.Range("E118").Value2 = 56
$ws.Range("G118").Value2 = 0
$ws.Range("H118").Value2 = 0
